$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2244.9092
$ws.Range("I40").Value = 2066.6667
$ws.Range("J40").Value = 2458.8
$ws.Range("K40").Value = 2066.6667
$ws.Range("L40").Value = 2458.8
$ws.Range("M40").Value = -1891.6667
$ws.Range("N40").Value = -2808.8
$ws.Range("H64").Value = 404381.7
$ws.Range("I64").Value = 836695.2
$ws.Range("J64").Value = 5323.077
$ws.Range("K64").Value = 836695.2
$ws.Range("L64").Value = 5323.077
$ws.Range("M64").Value = -836447.2
$ws.Range("N64").Value = -5819.077
$ws.Range("H67").Value = 404381.7
$ws.Range("I67").Value = 836695.2
$ws.Range("J67").Value = 5323.077
$ws.Range("K67").Value = 836695.2
$ws.Range("L67").Value = 5323.077
$ws.Range("M67").Value = -835837.2
$ws.Range("N67").Value = -7039.077
$ws.Range("H106").Value = 8626470
$ws.Range("I106").Value = 10194048
$ws.Range("J106").Value = 4788
$ws.Range("K106").Value = 10194048
$ws.Range("L106").Value = 4788
$ws.Range("M106").Value = -10193417
$ws.Range("N106").Value = -6050
$ws.Range("H133").Value = 40000
$ws.Range("J133").Value = 40000
$ws.Range("L133").Value = 40000
$ws.Range("N133").Value = -50120
$ws.Range("H136").Value = 41500
$ws.Range("J136").Value = 41500
$ws.Range("L136").Value = 41500
$ws.Range("N136").Value = -51700
$ws.Range("H137").Value = 28572796
$ws.Range("I137").Value = 37037980
$ws.Range("J137").Value = 2797
$ws.Range("K137").Value = 111113940
$ws.Range("L137").Value = 8391
$ws.Range("M137").Value = -111111390
$ws.Range("N137").Value = -13491
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4289.1455
$ws.Range("I32").Value = 2078.675
$ws.Range("J32").Value = 10183.733
$ws.Range("K32").Value = 2078.675
$ws.Range("L32").Value = 10183.733
$ws.Range("M32").Value = -1791.675
$ws.Range("N32").Value = -10757.733
$ws.Range("H61").Value = 2798.9768
$ws.Range("I61").Value = 2015.1072
$ws.Range("J61").Value = 4262.2
$ws.Range("K61").Value = 2015.1072
$ws.Range("L61").Value = 4262.2
$ws.Range("M61").Value = -1803.1072
$ws.Range("N61").Value = -4686.2
$ws.Range("H97").Value = 25649240
$ws.Range("I97").Value = 30312592
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 30312592
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -30312096
$ws.Range("N97").Value = -1792
$ws.Range("H102").Value = 2343.2
$ws.Range("I102").Value = 2269
$ws.Range("K102").Value = 2269
$ws.Range("M102").Value = -647
$ws.Range("H122").Value = 1194.4546
$ws.Range("I122").Value = 1073.1428
$ws.Range("J122").Value = 1406.75
$ws.Range("K122").Value = 3219.4284
$ws.Range("L122").Value = 4220.25
$ws.Range("M122").Value = -769.4284000000002
$ws.Range("N122").Value = -9120.25
$ws.Range("H132").Value = 2455.8604
$ws.Range("I132").Value = 1883.2812
$ws.Range("J132").Value = 4121.5454
$ws.Range("K132").Value = 5649.8436
$ws.Range("L132").Value = 12364.6362
$ws.Range("M132").Value = -3119.8436
$ws.Range("N132").Value = -17424.6362
$ws.Range("H136").Value = 2798.9768
$ws.Range("I136").Value = 2015.1072
$ws.Range("J136").Value = 4262.2
$ws.Range("K136").Value = 6045.321599999999
$ws.Range("L136").Value = 12786.6
$ws.Range("M136").Value = -3495.321599999999
$ws.Range("N136").Value = -17886.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1591.1428
$ws.Range("I94").Value = 1835.2727
$ws.Range("J94").Value = 696
$ws.Range("K94").Value = 1835.2727
$ws.Range("L94").Value = 696
$ws.Range("M94").Value = -1384.2727
$ws.Range("N94").Value = -1598
$ws.Range("H99").Value = 1417
$ws.Range("I99").Value = 1099.9
$ws.Range("K99").Value = 1099.9
$ws.Range("M99").Value = 398.0999999999999
$ws.Range("H105").Value = 3666.2666
$ws.Range("I105").Value = 3545.7273
$ws.Range("J105").Value = 3997.75
$ws.Range("K105").Value = 3545.7273
$ws.Range("L105").Value = 3997.75
$ws.Range("M105").Value = -1798.7273
$ws.Range("N105").Value = -7491.75
$ws.Range("H134").Value = 6299.154
$ws.Range("I134").Value = 8008
$ws.Range("J134").Value = 5786.5
$ws.Range("K134").Value = 24024
$ws.Range("L134").Value = 17359.5
$ws.Range("M134").Value = -21489
$ws.Range("N134").Value = -22429.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 922
$ws.Range("I16").Value = 537
$ws.Range("J16").Value = 1499.5
$ws.Range("K16").Value = 537
$ws.Range("L16").Value = 1499.5
$ws.Range("M16").Value = -250
$ws.Range("N16").Value = -2073.5
$ws.Range("H58").Value = 1834.9412
$ws.Range("I58").Value = 1326.0435
$ws.Range("J58").Value = 2899
$ws.Range("K58").Value = 1326.0435
$ws.Range("L58").Value = 2899
$ws.Range("M58").Value = -1123.0435
$ws.Range("N58").Value = -3305
$ws.Range("H62").Value = 13450.182
$ws.Range("I62").Value = 16393.666
$ws.Range("J62").Value = 7142.7144
$ws.Range("K62").Value = 16393.666
$ws.Range("L62").Value = 7142.7144
$ws.Range("M62").Value = -15769.666
$ws.Range("N62").Value = -8390.714400000001
$ws.Range("H65").Value = 13450.182
$ws.Range("I65").Value = 16393.666
$ws.Range("J65").Value = 7142.7144
$ws.Range("K65").Value = 81968.33
$ws.Range("L65").Value = 35713.572
$ws.Range("M65").Value = -78848.33
$ws.Range("N65").Value = -41953.572
$ws.Range("H113").Value = 922
$ws.Range("I113").Value = 537
$ws.Range("J113").Value = 1499.5
$ws.Range("K113").Value = 537
$ws.Range("L113").Value = 1499.5
$ws.Range("M113").Value = 1633
$ws.Range("N113").Value = -5839.5
$ws.Range("H122").Value = 2053.3076
$ws.Range("I122").Value = 899
$ws.Range("J122").Value = 3900.2
$ws.Range("K122").Value = 2697
$ws.Range("L122").Value = 11700.6
$ws.Range("M122").Value = -247
$ws.Range("N122").Value = -16600.6
$ws.Range("H136").Value = 1834.9412
$ws.Range("I136").Value = 1326.0435
$ws.Range("J136").Value = 2899
$ws.Range("K136").Value = 3978.1305
$ws.Range("L136").Value = 8697
$ws.Range("M136").Value = -1428.1305
$ws.Range("N136").Value = -13797
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1100.0952
$ws.Range("J5").Value = 2646
$ws.Range("L5").Value = 7938
$ws.Range("N5").Value = -8162
$ws.Range("H32").Value = 500
$ws.Range("J32").Value = 500
$ws.Range("L32").Value = 1500
$ws.Range("N32").Value = -2066
$ws.Range("H107").Value = 266.58334
$ws.Range("I107").Value = 211.17647
$ws.Range("K107").Value = 633.52941
$ws.Range("M107").Value = 1286.47059
$ws.Range("H113").Value = 13889670
$ws.Range("I113").Value = 411.5
$ws.Range("K113").Value = 1234.5
$ws.Range("M113").Value = 935.5
$ws.Range("H121").Value = 671.94116
$ws.Range("I121").Value = 262.85715
$ws.Range("J121").Value = 958.3
$ws.Range("K121").Value = 788.5714499999999
$ws.Range("L121").Value = 2874.9
$ws.Range("M121").Value = 521.4285500000001
$ws.Range("N121").Value = -5494.9
$ws.Range("H129").Value = 2467.9
$ws.Range("I129").Value = 1976.6666
$ws.Range("J129").Value = 2678.4285
$ws.Range("K129").Value = 5929.9998
$ws.Range("L129").Value = 8035.2855
$ws.Range("M129").Value = -929.9997999999996
$ws.Range("N129").Value = -18035.2855
$ws.Range("H131").Value = 1408.6143
$ws.Range("I131").Value = 578
$ws.Range("K131").Value = 1734
$ws.Range("M131").Value = 3306
$ws.Range("H133").Value = 7169.1665
$ws.Range("I133").Value = 1506
$ws.Range("J133").Value = 11214.286
$ws.Range("K133").Value = 4518
$ws.Range("L133").Value = 33642.858
$ws.Range("M133").Value = 542
$ws.Range("N133").Value = -43762.858
$ws.Range("H135").Value = 1100.0952
$ws.Range("J135").Value = 2646
$ws.Range("L135").Value = 23814
$ws.Range("N135").Value = -28884
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1389.125
$ws.Range("I102").Value = 1100
$ws.Range("J102").Value = 1430.4286
$ws.Range("K102").Value = 1100
$ws.Range("L102").Value = 1430.4286
$ws.Range("M102").Value = 522
$ws.Range("N102").Value = -4674.4286
$ws.Range("H132").Value = 2935.923
$ws.Range("I132").Value = 2502
$ws.Range("J132").Value = 4040.4546
$ws.Range("K132").Value = 7506
$ws.Range("L132").Value = 12121.3638
$ws.Range("M132").Value = -4976
$ws.Range("N132").Value = -17181.3638
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1229.8572
$ws.Range("I16").Value = 1592.3334
$ws.Range("J16").Value = 577.4
$ws.Range("K16").Value = 1592.3334
$ws.Range("L16").Value = 577.4
$ws.Range("M16").Value = -1422.3334
$ws.Range("N16").Value = -917.4
$ws.Range("H22").Value = 8921.076999999999
$ws.Range("J22").Value = 11238.4
$ws.Range("L22").Value = 11238.4
$ws.Range("N22").Value = -11828.4
$ws.Range("H27").Value = 8921.076999999999
$ws.Range("J27").Value = 11238.4
$ws.Range("L27").Value = 11238.4
$ws.Range("N27").Value = -11452.4
$ws.Range("H34").Value = 53333.332
$ws.Range("I34").Value = 53333.332
$ws.Range("K34").Value = 53333.332
$ws.Range("M34").Value = -53161.332
$ws.Range("H93").Value = 846.2222
$ws.Range("I93").Value = 768.6667
$ws.Range("J93").Value = 1234
$ws.Range("K93").Value = 768.6667
$ws.Range("L93").Value = 1234
$ws.Range("M93").Value = 479.3333
$ws.Range("N93").Value = -3730
$ws.Range("H100").Value = 2439.2122
$ws.Range("I100").Value = 1942.4286
$ws.Range("J100").Value = 2805.2632
$ws.Range("K100").Value = 1942.4286
$ws.Range("L100").Value = 2805.2632
$ws.Range("M100").Value = -1401.4286
$ws.Range("N100").Value = -3887.2632
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10419041
$ws.Range("I132").Value = 12822411
$ws.Range("J132").Value = 4436.5557
$ws.Range("K132").Value = 38467233
$ws.Range("L132").Value = 13309.6671
$ws.Range("M132").Value = -38464703
$ws.Range("N132").Value = -18369.6671
